$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("C2").Value = "Hartmut"
# Card number is long enough that Excel would otherwise render it in
# scientific notation as a plain number, so force text storage first.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 21.11.2023"

# Row 6
$ws.Range("B6").Value = "24.11."
$ws.Range("C6").Value = "25.11."
$ws.Range("D6").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E6").Value = "58,63-"

# Row 7
$ws.Range("B7").Value = "27.11."
$ws.Range("C7").Value = "28.11."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-28252062"
$ws.Range("E7").Value = "54,58-"

# Row 8
$ws.Range("B8").Value = "01.12."
$ws.Range("C8").Value = "02.12."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,16-"

# Row 9
$ws.Range("B9").Value = "05.12."
$ws.Range("C9").Value = "06.12."
$ws.Range("D9").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E9").Value = "62,42-"

# Row 10
$ws.Range("B10").Value = "06.12."
$ws.Range("C10").Value = "07.12."
$ws.Range("D10").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E10").Value = "69,94-"

# Row 11 - newly populated, copy style from row 10 so formatting matches other data rows
$ws.Range("B10:E10").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B11").Value = "09.12."
$ws.Range("C11").Value = "10.12."
$ws.Range("D11").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E11").Value = "44,02-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 14.12.2023"
$ws.Range("E12").Value = "314,75-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 22.12.2023"
